# Add 2022-Q4 data:
#  1. Insert a new worksheet named "2022-Q4" right after "总计" (i.e. before the
#     sheet currently named "2022-Q3"), cloned from the "2022-Q3" sheet so it
#     inherits the same column layout/styles, then overwrite its data.
#  2. Prepend a new "2022-Q4" row to the "总计" summary sheet (shifting the
#     existing quarters down by one row) and fill in the new totals.

$wb = $excel.ActiveWorkbook

# Remember the workbook's current active sheet so we can restore the
# selection/active-tab state once we're done (copying a sheet makes the new
# copy the active sheet as a side effect).
$originalActiveSheetName = $wb.ActiveSheet.Name

# --- 1. Clone the "2022-Q3" sheet to create "2022-Q4" ----------------------

$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet)

$q4Sheet = $wb.Worksheets.Item("2022-Q3 (2)")
$q4Sheet.Name = "2022-Q4"

# Clear out the cloned "2022-Q3" rows (2 data rows) before writing the new
# 2022-Q4 figures (6 data rows).
$q4Sheet.Range("A2:H3").ClearContents()

$q4Data = @(
    @("003751", "万家瑞隆混合A",          "14.93", "90.82", "2.47", "0.3688", 8),
    @("015384", "万家瑞隆混合C",          "3.16",  "90.82", "2.47", "0.0781", 8),
    @("013326", "万家景气驱动混合A",      "2.03",  "92.56", "2.66", "0.0540", 9),
    @("008491", "万家周期优势企业混合A",  "0.62",  "91.08", "2.73", "0.0169", 8),
    @("013327", "万家景气驱动混合C",      "0.24",  "92.56", "2.66", "0.0064", 9),
    @("008492", "万家周期优势企业混合C",  "0.15",  "91.08", "2.73", "0.0041", 8)
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $row = $i + 2
    $vals = $q4Data[$i]

    if ($i -gt 1) {
        # Rows 4-7 don't exist yet on the cloned sheet (it only had 2 data
        # rows) - stamp the row style from row 3 (A3 has the "index column"
        # style) before filling it in.
        $q4Sheet.Range("A3").Copy($q4Sheet.Range("A" + $row))
    }

    $q4Sheet.Cells.Item($row, 1).Value = $i
    $q4Sheet.Cells.Item($row, 2).Value = "'" + $vals[0]
    $q4Sheet.Cells.Item($row, 3).Value = $vals[1]
    $q4Sheet.Cells.Item($row, 4).Value = "'" + $vals[2]
    $q4Sheet.Cells.Item($row, 5).Value = "'" + $vals[3]
    $q4Sheet.Cells.Item($row, 6).Value = "'" + $vals[4]
    $q4Sheet.Cells.Item($row, 7).Value = "'" + $vals[5]
    $q4Sheet.Cells.Item($row, 8).Value = $vals[6]
}

# --- 2. Shift "总计" rows down and insert the new 2022-Q4 totals -----------
# (Values are written explicitly rather than copied from the existing rows,
# since the known end-state is fully determined by the source diff.)

$totalSheet = $wb.Worksheets.Item("总计")

# Stamp the new row 7 (index column) from row 6's style before filling it in.
$totalSheet.Range("A6").Copy($totalSheet.Range("A7"))

$totalRows = @(
    @(0, "2022-Q4", 6, 0.53),
    @(1, "2022-Q3", 2, 1.39),
    @(2, "2022-Q2", 3, 1.8),
    @(3, "2022-Q1", 4, 4.6),
    @(4, "2021-Q4", 1, 0.5600000000000001),
    @(5, "2021-Q3", 1, 0.44)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $row = $i + 2
    $vals = $totalRows[$i]
    $totalSheet.Cells.Item($row, 1).Value = $vals[0]
    $totalSheet.Cells.Item($row, 2).Value = $vals[1]
    $totalSheet.Cells.Item($row, 3).Value = $vals[2]
    $totalSheet.Cells.Item($row, 4).Value = $vals[3]
}

# --- Restore the original active sheet/selection ---------------------------

$wb.Worksheets.Item($originalActiveSheetName).Activate()
